$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.657.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.591.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.60%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.45%  '
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -2.71%  '
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("E10").Value = '  -3.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0836'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.813.84'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.590.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.47%  '
$ws.Range("E14").Value = '  -2.20%  '
$ws.Range("E15").Value = '  -3.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.640.39'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.08%  '
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '207.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.24%  '
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.84%  '
$ws.Range("E22").Value = '  -2.97%  '
$ws.Range("E23").Value = '  -3.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("E28").Value = '  -3.47%  '
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0503'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.82%  '
$ws.Range("E31").Value = '  -2.09%  '
$ws.Range("E32").Value = '  -4.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.665'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +22.53%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.95%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.322.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.41'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.93%  '
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.828'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.46%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.39'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.31%  '
$ws.Range("E42").Value = '  -1.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.726.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("E46").Value = '  -0.89%  '
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.840'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0974'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.63%  '
